$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241001-093545-"

# Update the snapshot date (column G) for every data row from 45565 (2024-09-30) to 45566 (2024-10-01)
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45566
}

# Update refreshed balance values (columns E and H mirror each other) for the rows whose balances changed
$balanceUpdates = @{
    5 = 1448
    6 = 1059.1400000000001
    8 = 155.22999999999999
    15 = 1570.97
    17 = 1174.93
    36 = 0.02
    43 = 822.38
    51 = 1017.18
    52 = 924.47
    57 = 3973.31
    58 = 33.51
    60 = 1130.1400000000001
    97 = 1256.82
    99 = 1615.76
    101 = 7170.3
    102 = 773.62
    103 = 1073.05
    104 = 230.37
    105 = 642.86
    107 = 1178.46
    108 = 439.25
    109 = 30639.39
    110 = 820.72
    112 = 0.42
    113 = 3.22
    118 = 916.88
    132 = 1015.86
    138 = 1661.45
    143 = 1943.09
    148 = -6283.58
    158 = 85.93
    161 = 289.57
    165 = 1042.07
    173 = 2216.06
    224 = 605.29999999999995
    230 = 848.24
    232 = 40.26
    235 = 590.04999999999995
    249 = 932.35
    255 = 57457.04
    264 = 2951.63
    265 = 1965.65
    270 = 32.909999999999997
    271 = 1322.69
    273 = 1627.21
    274 = 0
}

foreach ($row in $balanceUpdates.Keys) {
    $value = $balanceUpdates[$row]
    $ws.Cells.Item($row, 5).Value = $value   # column E
    $ws.Cells.Item($row, 8).Value = $value   # column H
}

Write-Host "Applied balance refresh updates"
